# Fill in the Região (A) column for the rows that were missing it.
# Shared strings already contain "Sergipe", "Nordeste" and "Brasil",
# so setting .Value re-uses those entries instead of creating new ones.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($r in 3..6) {
    $ws.Cells.Item($r, 1).Value = "Sergipe"
}

foreach ($r in 8..11) {
    $ws.Cells.Item($r, 1).Value = "Brasil"
}

foreach ($r in 13..16) {
    $ws.Cells.Item($r, 1).Value = "Nordeste"
}

# Update the active selection to match the edited region.
$ws.Range("A12:A16").Select()
